$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Cells.Item(75, 1).Value = 7
$ws.Cells.Item(75, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(75, 3).Value = 'Ñuble'
$ws.Cells.Item(75, 4).Value = 44586
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 100112024
$ws.Cells.Item(75, 7).Value = 'Choclo'
$ws.Cells.Item(75, 8).Value = 'Choclero'
$ws.Cells.Item(75, 9).Value = 'Primera'
$ws.Cells.Item(75, 10).Value = 26000
$ws.Cells.Item(75, 11).Value = 200
$ws.Cells.Item(75, 12).Value = 250
$ws.Cells.Item(75, 13).Value = 225
$ws.Cells.Item(75, 14).Value = '$/unidad'
$ws.Cells.Item(75, 15).Value = 'Región del Maule'
$ws.Cells.Item(75, 16).Value = 225
$ws.Cells.Item(75, 17).Value = 1
$ws.Cells.Item(75, 18).Value = 'Hortaliza'

# Row 76
$ws.Cells.Item(76, 1).Value = 7
$ws.Cells.Item(76, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(76, 3).Value = 'Ñuble'
$ws.Cells.Item(76, 4).Value = 44586
$ws.Cells.Item(76, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(76, 5).Value = 16
$ws.Cells.Item(76, 6).Value = 100112024
$ws.Cells.Item(76, 7).Value = 'Choclo'
$ws.Cells.Item(76, 8).Value = 'Choclero'
$ws.Cells.Item(76, 9).Value = 'Segunda'
$ws.Cells.Item(76, 10).Value = 20000
$ws.Cells.Item(76, 11).Value = 100
$ws.Cells.Item(76, 12).Value = 150
$ws.Cells.Item(76, 13).Value = 125
$ws.Cells.Item(76, 14).Value = '$/unidad'
$ws.Cells.Item(76, 15).Value = 'Región del Maule'
$ws.Cells.Item(76, 16).Value = 125
$ws.Cells.Item(76, 17).Value = 1
$ws.Cells.Item(76, 18).Value = 'Hortaliza'

# Row 77
$ws.Cells.Item(77, 1).Value = 7
$ws.Cells.Item(77, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(77, 3).Value = 'Ñuble'
$ws.Cells.Item(77, 4).Value = 44558
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 5).Value = 16
$ws.Cells.Item(77, 6).Value = 100112024
$ws.Cells.Item(77, 7).Value = 'Choclo'
$ws.Cells.Item(77, 8).Value = 'Choclero'
$ws.Cells.Item(77, 9).Value = 'Primera'
$ws.Cells.Item(77, 10).Value = 12000
$ws.Cells.Item(77, 11).Value = 300
$ws.Cells.Item(77, 12).Value = 350
$ws.Cells.Item(77, 13).Value = 325
$ws.Cells.Item(77, 14).Value = '$/unidad'
$ws.Cells.Item(77, 15).Value = 'Región del Maule'
$ws.Cells.Item(77, 16).Value = 325
$ws.Cells.Item(77, 17).Value = 1
$ws.Cells.Item(77, 18).Value = 'Hortaliza'

# Row 78
$ws.Cells.Item(78, 1).Value = 7
$ws.Cells.Item(78, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(78, 3).Value = 'Ñuble'
$ws.Cells.Item(78, 4).Value = 44558
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 16
$ws.Cells.Item(78, 6).Value = 100112024
$ws.Cells.Item(78, 7).Value = 'Choclo'
$ws.Cells.Item(78, 8).Value = 'Choclero'
$ws.Cells.Item(78, 9).Value = 'Segunda'
$ws.Cells.Item(78, 10).Value = 4000
$ws.Cells.Item(78, 11).Value = 200
$ws.Cells.Item(78, 12).Value = 250
$ws.Cells.Item(78, 13).Value = 225
$ws.Cells.Item(78, 14).Value = '$/unidad'
$ws.Cells.Item(78, 15).Value = 'Región del Maule'
$ws.Cells.Item(78, 16).Value = 225
$ws.Cells.Item(78, 17).Value = 1
$ws.Cells.Item(78, 18).Value = 'Hortaliza'

# Row 79
$ws.Cells.Item(79, 1).Value = 7
$ws.Cells.Item(79, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(79, 3).Value = 'Ñuble'
$ws.Cells.Item(79, 4).Value = 44524
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 5).Value = 16
$ws.Cells.Item(79, 6).Value = 100112024
$ws.Cells.Item(79, 7).Value = 'Choclo'
$ws.Cells.Item(79, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(79, 9).Value = 'Primera'
$ws.Cells.Item(79, 10).Value = 60
$ws.Cells.Item(79, 11).Value = 16000
$ws.Cells.Item(79, 12).Value = 17000
$ws.Cells.Item(79, 13).Value = 16500
$ws.Cells.Item(79, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(79, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(79, 16).Value = 275
$ws.Cells.Item(79, 17).Value = 60
$ws.Cells.Item(79, 18).Value = 'Hortaliza'

# Row 80
$ws.Cells.Item(80, 1).Value = 7
$ws.Cells.Item(80, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(80, 3).Value = 'Ñuble'
$ws.Cells.Item(80, 4).Value = 44530
$ws.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80, 5).Value = 16
$ws.Cells.Item(80, 6).Value = 100112024
$ws.Cells.Item(80, 7).Value = 'Choclo'
$ws.Cells.Item(80, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(80, 9).Value = 'Primera'
$ws.Cells.Item(80, 10).Value = 60
$ws.Cells.Item(80, 11).Value = 16000
$ws.Cells.Item(80, 12).Value = 17000
$ws.Cells.Item(80, 13).Value = 16500
$ws.Cells.Item(80, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(80, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(80, 16).Value = 275
$ws.Cells.Item(80, 17).Value = 60
$ws.Cells.Item(80, 18).Value = 'Hortaliza'

# Row 81
$ws.Cells.Item(81, 1).Value = 7
$ws.Cells.Item(81, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(81, 3).Value = 'Ñuble'
$ws.Cells.Item(81, 4).Value = 44267
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 5).Value = 16
$ws.Cells.Item(81, 6).Value = 100112024
$ws.Cells.Item(81, 7).Value = 'Choclo'
$ws.Cells.Item(81, 8).Value = 'Choclero'
$ws.Cells.Item(81, 9).Value = 'Primera'
$ws.Cells.Item(81, 10).Value = 53000
$ws.Cells.Item(81, 11).Value = 250
$ws.Cells.Item(81, 12).Value = 270
$ws.Cells.Item(81, 13).Value = 257
$ws.Cells.Item(81, 14).Value = '$/unidad'
$ws.Cells.Item(81, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(81, 16).Value = 257
$ws.Cells.Item(81, 17).Value = 1
$ws.Cells.Item(81, 18).Value = 'Hortaliza'

# Row 82
$ws.Cells.Item(82, 1).Value = 7
$ws.Cells.Item(82, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(82, 3).Value = 'Ñuble'
$ws.Cells.Item(82, 4).Value = 44267
$ws.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82, 5).Value = 16
$ws.Cells.Item(82, 6).Value = 100112024
$ws.Cells.Item(82, 7).Value = 'Choclo'
$ws.Cells.Item(82, 8).Value = 'Choclero'
$ws.Cells.Item(82, 9).Value = 'Segunda'
$ws.Cells.Item(82, 10).Value = 15000
$ws.Cells.Item(82, 11).Value = 200
$ws.Cells.Item(82, 12).Value = 200
$ws.Cells.Item(82, 13).Value = 200
$ws.Cells.Item(82, 14).Value = '$/unidad'
$ws.Cells.Item(82, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(82, 16).Value = 200
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = 'Hortaliza'

# Row 83
$ws.Cells.Item(83, 1).Value = 7
$ws.Cells.Item(83, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(83, 3).Value = 'Ñuble'
$ws.Cells.Item(83, 4).Value = 44235
$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83, 5).Value = 16
$ws.Cells.Item(83, 6).Value = 100112024
$ws.Cells.Item(83, 7).Value = 'Choclo'
$ws.Cells.Item(83, 8).Value = 'Choclero'
$ws.Cells.Item(83, 9).Value = 'Primera'
$ws.Cells.Item(83, 10).Value = 16000
$ws.Cells.Item(83, 11).Value = 240
$ws.Cells.Item(83, 12).Value = 250
$ws.Cells.Item(83, 13).Value = 245
$ws.Cells.Item(83, 14).Value = '$/unidad'
$ws.Cells.Item(83, 15).Value = 'Región del Maule'
$ws.Cells.Item(83, 16).Value = 245
$ws.Cells.Item(83, 17).Value = 1
$ws.Cells.Item(83, 18).Value = 'Hortaliza'

# Row 84
$ws.Cells.Item(84, 1).Value = 7
$ws.Cells.Item(84, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(84, 3).Value = 'Ñuble'
$ws.Cells.Item(84, 4).Value = 44580
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 5).Value = 16
$ws.Cells.Item(84, 6).Value = 100112024
$ws.Cells.Item(84, 7).Value = 'Choclo'
$ws.Cells.Item(84, 8).Value = 'Choclero'
$ws.Cells.Item(84, 9).Value = 'Primera'
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 11).Value = 250
$ws.Cells.Item(84, 12).Value = 300
$ws.Cells.Item(84, 13).Value = 275
$ws.Cells.Item(84, 14).Value = '$/unidad'
$ws.Cells.Item(84, 15).Value = 'Región del Maule'
$ws.Cells.Item(84, 16).Value = 275
$ws.Cells.Item(84, 17).Value = 1
$ws.Cells.Item(84, 18).Value = 'Hortaliza'

# Row 85
$ws.Cells.Item(85, 1).Value = 7
$ws.Cells.Item(85, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(85, 3).Value = 'Ñuble'
$ws.Cells.Item(85, 4).Value = 44580
$ws.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85, 5).Value = 16
$ws.Cells.Item(85, 6).Value = 100112024
$ws.Cells.Item(85, 7).Value = 'Choclo'
$ws.Cells.Item(85, 8).Value = 'Choclero'
$ws.Cells.Item(85, 9).Value = 'Segunda'
$ws.Cells.Item(85, 10).Value = 10000
$ws.Cells.Item(85, 11).Value = 200
$ws.Cells.Item(85, 12).Value = 200
$ws.Cells.Item(85, 13).Value = 200
$ws.Cells.Item(85, 14).Value = '$/unidad'
$ws.Cells.Item(85, 15).Value = 'Región del Maule'
$ws.Cells.Item(85, 16).Value = 200
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = 'Hortaliza'

# Row 86
$ws.Cells.Item(86, 1).Value = 7
$ws.Cells.Item(86, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(86, 3).Value = 'Ñuble'
$ws.Cells.Item(86, 4).Value = 44309
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 5).Value = 16
$ws.Cells.Item(86, 6).Value = 100112024
$ws.Cells.Item(86, 7).Value = 'Choclo'
$ws.Cells.Item(86, 8).Value = 'Choclero'
$ws.Cells.Item(86, 9).Value = 'Primera'
$ws.Cells.Item(86, 10).Value = 6000
$ws.Cells.Item(86, 11).Value = 230
$ws.Cells.Item(86, 12).Value = 250
$ws.Cells.Item(86, 13).Value = 240
$ws.Cells.Item(86, 14).Value = '$/unidad'
$ws.Cells.Item(86, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(86, 16).Value = 240
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = 'Hortaliza'

# Row 87
$ws.Cells.Item(87, 1).Value = 7
$ws.Cells.Item(87, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(87, 3).Value = 'Ñuble'
$ws.Cells.Item(87, 4).Value = 44313
$ws.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(87, 5).Value = 16
$ws.Cells.Item(87, 6).Value = 100112024
$ws.Cells.Item(87, 7).Value = 'Choclo'
$ws.Cells.Item(87, 8).Value = 'Choclero'
$ws.Cells.Item(87, 9).Value = 'Primera'
$ws.Cells.Item(87, 10).Value = 3000
$ws.Cells.Item(87, 11).Value = 200
$ws.Cells.Item(87, 12).Value = 200
$ws.Cells.Item(87, 13).Value = 200
$ws.Cells.Item(87, 14).Value = '$/unidad'
$ws.Cells.Item(87, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(87, 16).Value = 200
$ws.Cells.Item(87, 17).Value = 1
$ws.Cells.Item(87, 18).Value = 'Hortaliza'

# Row 88
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(88, 3).Value = 'Ñuble'
$ws.Cells.Item(88, 4).Value = 44242
$ws.Cells.Item(88, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = 100112024
$ws.Cells.Item(88, 7).Value = 'Choclo'
$ws.Cells.Item(88, 8).Value = 'Choclero'
$ws.Cells.Item(88, 9).Value = 'Primera'
$ws.Cells.Item(88, 10).Value = 16000
$ws.Cells.Item(88, 11).Value = 300
$ws.Cells.Item(88, 12).Value = 300
$ws.Cells.Item(88, 13).Value = 300
$ws.Cells.Item(88, 14).Value = '$/unidad'
$ws.Cells.Item(88, 15).Value = 'Región del Maule'
$ws.Cells.Item(88, 16).Value = 300
$ws.Cells.Item(88, 17).Value = 1
$ws.Cells.Item(88, 18).Value = 'Hortaliza'

# Row 89
$ws.Cells.Item(89, 1).Value = 7
$ws.Cells.Item(89, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(89, 3).Value = 'Ñuble'
$ws.Cells.Item(89, 4).Value = 44572
$ws.Cells.Item(89, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(89, 5).Value = 16
$ws.Cells.Item(89, 6).Value = 100112024
$ws.Cells.Item(89, 7).Value = 'Choclo'
$ws.Cells.Item(89, 8).Value = 'Choclero'
$ws.Cells.Item(89, 9).Value = 'Primera'
$ws.Cells.Item(89, 10).Value = 12000
$ws.Cells.Item(89, 11).Value = 200
$ws.Cells.Item(89, 12).Value = 250
$ws.Cells.Item(89, 13).Value = 225
$ws.Cells.Item(89, 14).Value = '$/unidad'
$ws.Cells.Item(89, 15).Value = 'Región del Maule'
$ws.Cells.Item(89, 16).Value = 225
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = 'Hortaliza'

# Row 90
$ws.Cells.Item(90, 1).Value = 7
$ws.Cells.Item(90, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(90, 3).Value = 'Ñuble'
$ws.Cells.Item(90, 4).Value = 44572
$ws.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(90, 5).Value = 16
$ws.Cells.Item(90, 6).Value = 100112024
$ws.Cells.Item(90, 7).Value = 'Choclo'
$ws.Cells.Item(90, 8).Value = 'Choclero'
$ws.Cells.Item(90, 9).Value = 'Segunda'
$ws.Cells.Item(90, 10).Value = 10000
$ws.Cells.Item(90, 11).Value = 100
$ws.Cells.Item(90, 12).Value = 150
$ws.Cells.Item(90, 13).Value = 125
$ws.Cells.Item(90, 14).Value = '$/unidad'
$ws.Cells.Item(90, 15).Value = 'Región del Maule'
$ws.Cells.Item(90, 16).Value = 125
$ws.Cells.Item(90, 17).Value = 1
$ws.Cells.Item(90, 18).Value = 'Hortaliza'

# Row 91
$ws.Cells.Item(91, 1).Value = 7
$ws.Cells.Item(91, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(91, 3).Value = 'Ñuble'
$ws.Cells.Item(91, 4).Value = 44257
$ws.Cells.Item(91, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(91, 5).Value = 16
$ws.Cells.Item(91, 6).Value = 100112024
$ws.Cells.Item(91, 7).Value = 'Choclo'
$ws.Cells.Item(91, 8).Value = 'Choclero'
$ws.Cells.Item(91, 9).Value = 'Primera'
$ws.Cells.Item(91, 10).Value = 40000
$ws.Cells.Item(91, 11).Value = 260
$ws.Cells.Item(91, 12).Value = 280
$ws.Cells.Item(91, 13).Value = 271
$ws.Cells.Item(91, 14).Value = '$/unidad'
$ws.Cells.Item(91, 15).Value = 'Región del Maule'
$ws.Cells.Item(91, 16).Value = 271
$ws.Cells.Item(91, 17).Value = 1
$ws.Cells.Item(91, 18).Value = 'Hortaliza'

# Row 92
$ws.Cells.Item(92, 1).Value = 7
$ws.Cells.Item(92, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(92, 3).Value = 'Ñuble'
$ws.Cells.Item(92, 4).Value = 44552
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 16
$ws.Cells.Item(92, 6).Value = 100112024
$ws.Cells.Item(92, 7).Value = 'Choclo'
$ws.Cells.Item(92, 8).Value = 'Choclero'
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 10000
$ws.Cells.Item(92, 11).Value = 300
$ws.Cells.Item(92, 12).Value = 350
$ws.Cells.Item(92, 13).Value = 325
$ws.Cells.Item(92, 14).Value = '$/unidad'
$ws.Cells.Item(92, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(92, 16).Value = 325
$ws.Cells.Item(92, 17).Value = 1
$ws.Cells.Item(92, 18).Value = 'Hortaliza'

# Row 93
$ws.Cells.Item(93, 1).Value = 7
$ws.Cells.Item(93, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(93, 3).Value = 'Ñuble'
$ws.Cells.Item(93, 4).Value = 44552
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(93, 5).Value = 16
$ws.Cells.Item(93, 6).Value = 100112024
$ws.Cells.Item(93, 7).Value = 'Choclo'
$ws.Cells.Item(93, 8).Value = 'Choclero'
$ws.Cells.Item(93, 9).Value = 'Segunda'
$ws.Cells.Item(93, 10).Value = 6000
$ws.Cells.Item(93, 11).Value = 200
$ws.Cells.Item(93, 12).Value = 250
$ws.Cells.Item(93, 13).Value = 225
$ws.Cells.Item(93, 14).Value = '$/unidad'
$ws.Cells.Item(93, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(93, 16).Value = 225
$ws.Cells.Item(93, 17).Value = 1
$ws.Cells.Item(93, 18).Value = 'Hortaliza'

# Row 94
$ws.Cells.Item(94, 1).Value = 7
$ws.Cells.Item(94, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(94, 3).Value = 'Ñuble'
$ws.Cells.Item(94, 4).Value = 44322
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 5).Value = 16
$ws.Cells.Item(94, 6).Value = 100112024
$ws.Cells.Item(94, 7).Value = 'Choclo'
$ws.Cells.Item(94, 8).Value = 'Choclero'
$ws.Cells.Item(94, 9).Value = 'Primera'
$ws.Cells.Item(94, 10).Value = 1400
$ws.Cells.Item(94, 11).Value = 180
$ws.Cells.Item(94, 12).Value = 200
$ws.Cells.Item(94, 13).Value = 190
$ws.Cells.Item(94, 14).Value = '$/unidad'
$ws.Cells.Item(94, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(94, 16).Value = 190
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = 'Hortaliza'

# Row 95
$ws.Cells.Item(95, 1).Value = 7
$ws.Cells.Item(95, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(95, 3).Value = 'Ñuble'
$ws.Cells.Item(95, 4).Value = 44246
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 5).Value = 16
$ws.Cells.Item(95, 6).Value = 100112024
$ws.Cells.Item(95, 7).Value = 'Choclo'
$ws.Cells.Item(95, 8).Value = 'Choclero'
$ws.Cells.Item(95, 9).Value = 'Primera'
$ws.Cells.Item(95, 10).Value = 16000
$ws.Cells.Item(95, 11).Value = 340
$ws.Cells.Item(95, 12).Value = 350
$ws.Cells.Item(95, 13).Value = 345
$ws.Cells.Item(95, 14).Value = '$/unidad'
$ws.Cells.Item(95, 15).Value = 'Región del Maule'
$ws.Cells.Item(95, 16).Value = 345
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = 'Hortaliza'

# Row 96
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(96, 3).Value = 'Ñuble'
$ws.Cells.Item(96, 4).Value = 44218
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112024
$ws.Cells.Item(96, 7).Value = 'Choclo'
$ws.Cells.Item(96, 8).Value = 'Choclero'
$ws.Cells.Item(96, 9).Value = 'Primera'
$ws.Cells.Item(96, 10).Value = 75000
$ws.Cells.Item(96, 11).Value = 270
$ws.Cells.Item(96, 12).Value = 300
$ws.Cells.Item(96, 13).Value = 280
$ws.Cells.Item(96, 14).Value = '$/unidad'
$ws.Cells.Item(96, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(96, 16).Value = 280
$ws.Cells.Item(96, 17).Value = 1
$ws.Cells.Item(96, 18).Value = 'Hortaliza'

# Row 97
$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(97, 3).Value = 'Ñuble'
$ws.Cells.Item(97, 4).Value = 44218
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = 100112024
$ws.Cells.Item(97, 7).Value = 'Choclo'
$ws.Cells.Item(97, 8).Value = 'Choclero'
$ws.Cells.Item(97, 9).Value = 'Segunda'
$ws.Cells.Item(97, 10).Value = 27000
$ws.Cells.Item(97, 11).Value = 200
$ws.Cells.Item(97, 12).Value = 230
$ws.Cells.Item(97, 13).Value = 217
$ws.Cells.Item(97, 14).Value = '$/unidad'
$ws.Cells.Item(97, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(97, 16).Value = 217
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = 'Hortaliza'

# Row 98
$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(98, 3).Value = 'Ñuble'
$ws.Cells.Item(98, 4).Value = 44211
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = 100112024
$ws.Cells.Item(98, 7).Value = 'Choclo'
$ws.Cells.Item(98, 8).Value = 'Choclero'
$ws.Cells.Item(98, 9).Value = 'Primera'
$ws.Cells.Item(98, 10).Value = 60000
$ws.Cells.Item(98, 11).Value = 280
$ws.Cells.Item(98, 12).Value = 300
$ws.Cells.Item(98, 13).Value = 292
$ws.Cells.Item(98, 14).Value = '$/unidad'
$ws.Cells.Item(98, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(98, 16).Value = 292
$ws.Cells.Item(98, 17).Value = 1
$ws.Cells.Item(98, 18).Value = 'Hortaliza'

# Row 99
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(99, 3).Value = 'Ñuble'
$ws.Cells.Item(99, 4).Value = 44211
$ws.Cells.Item(99, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112024
$ws.Cells.Item(99, 7).Value = 'Choclo'
$ws.Cells.Item(99, 8).Value = 'Choclero'
$ws.Cells.Item(99, 9).Value = 'Segunda'
$ws.Cells.Item(99, 10).Value = 30000
$ws.Cells.Item(99, 11).Value = 200
$ws.Cells.Item(99, 12).Value = 230
$ws.Cells.Item(99, 13).Value = 218
$ws.Cells.Item(99, 14).Value = '$/unidad'
$ws.Cells.Item(99, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(99, 16).Value = 218
$ws.Cells.Item(99, 17).Value = 1
$ws.Cells.Item(99, 18).Value = 'Hortaliza'

# Row 100
$ws.Cells.Item(100, 1).Value = 7
$ws.Cells.Item(100, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(100, 3).Value = 'Ñuble'
$ws.Cells.Item(100, 4).Value = 44540
$ws.Cells.Item(100, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(100, 5).Value = 16
$ws.Cells.Item(100, 6).Value = 100112024
$ws.Cells.Item(100, 7).Value = 'Choclo'
$ws.Cells.Item(100, 8).Value = 'Choclero'
$ws.Cells.Item(100, 9).Value = 'Primera'
$ws.Cells.Item(100, 10).Value = 600
$ws.Cells.Item(100, 11).Value = 450
$ws.Cells.Item(100, 12).Value = 500
$ws.Cells.Item(100, 13).Value = 475
$ws.Cells.Item(100, 14).Value = '$/unidad'
$ws.Cells.Item(100, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(100, 16).Value = 475
$ws.Cells.Item(100, 17).Value = 1
$ws.Cells.Item(100, 18).Value = 'Hortaliza'

# Row 101
$ws.Cells.Item(101, 1).Value = 7
$ws.Cells.Item(101, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(101, 3).Value = 'Ñuble'
$ws.Cells.Item(101, 4).Value = 44193
$ws.Cells.Item(101, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(101, 5).Value = 16
$ws.Cells.Item(101, 6).Value = 100112024
$ws.Cells.Item(101, 7).Value = 'Choclo'
$ws.Cells.Item(101, 8).Value = 'Choclero'
$ws.Cells.Item(101, 9).Value = 'Primera'
$ws.Cells.Item(101, 10).Value = 6000
$ws.Cells.Item(101, 11).Value = 400
$ws.Cells.Item(101, 12).Value = 420
$ws.Cells.Item(101, 13).Value = 410
$ws.Cells.Item(101, 14).Value = '$/unidad'
$ws.Cells.Item(101, 15).Value = 'Región del Maule'
$ws.Cells.Item(101, 16).Value = 410
$ws.Cells.Item(101, 17).Value = 1
$ws.Cells.Item(101, 18).Value = 'Hortaliza'

# Row 102
$ws.Cells.Item(102, 1).Value = 7
$ws.Cells.Item(102, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(102, 3).Value = 'Ñuble'
$ws.Cells.Item(102, 4).Value = 44200
$ws.Cells.Item(102, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(102, 5).Value = 16
$ws.Cells.Item(102, 6).Value = 100112024
$ws.Cells.Item(102, 7).Value = 'Choclo'
$ws.Cells.Item(102, 8).Value = 'Choclero'
$ws.Cells.Item(102, 9).Value = 'Primera'
$ws.Cells.Item(102, 10).Value = 40000
$ws.Cells.Item(102, 11).Value = 300
$ws.Cells.Item(102, 12).Value = 350
$ws.Cells.Item(102, 13).Value = 319
$ws.Cells.Item(102, 14).Value = '$/unidad'
$ws.Cells.Item(102, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(102, 16).Value = 319
$ws.Cells.Item(102, 17).Value = 1
$ws.Cells.Item(102, 18).Value = 'Hortaliza'

# Row 103
$ws.Cells.Item(103, 1).Value = 7
$ws.Cells.Item(103, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(103, 3).Value = 'Ñuble'
$ws.Cells.Item(103, 4).Value = 44200
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 5).Value = 16
$ws.Cells.Item(103, 6).Value = 100112024
$ws.Cells.Item(103, 7).Value = 'Choclo'
$ws.Cells.Item(103, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(103, 9).Value = 'Primera'
$ws.Cells.Item(103, 10).Value = 25000
$ws.Cells.Item(103, 11).Value = 280
$ws.Cells.Item(103, 12).Value = 300
$ws.Cells.Item(103, 13).Value = 290
$ws.Cells.Item(103, 14).Value = '$/unidad'
$ws.Cells.Item(103, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(103, 16).Value = 290
$ws.Cells.Item(103, 17).Value = 1
$ws.Cells.Item(103, 18).Value = 'Hortaliza'

# Row 104
$ws.Cells.Item(104, 1).Value = 7
$ws.Cells.Item(104, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(104, 3).Value = 'Ñuble'
$ws.Cells.Item(104, 4).Value = 44160
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 5).Value = 16
$ws.Cells.Item(104, 6).Value = 100112024
$ws.Cells.Item(104, 7).Value = 'Choclo'
$ws.Cells.Item(104, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(104, 9).Value = 'Primera'
$ws.Cells.Item(104, 10).Value = 45
$ws.Cells.Item(104, 11).Value = 24000
$ws.Cells.Item(104, 12).Value = 25000
$ws.Cells.Item(104, 13).Value = 24444
$ws.Cells.Item(104, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(104, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(104, 16).Value = 349
$ws.Cells.Item(104, 17).Value = 70
$ws.Cells.Item(104, 18).Value = 'Hortaliza'

# Row 105
$ws.Cells.Item(105, 1).Value = 7
$ws.Cells.Item(105, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(105, 3).Value = 'Ñuble'
$ws.Cells.Item(105, 4).Value = 44203
$ws.Cells.Item(105, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(105, 5).Value = 16
$ws.Cells.Item(105, 6).Value = 100112024
$ws.Cells.Item(105, 7).Value = 'Choclo'
$ws.Cells.Item(105, 8).Value = 'Choclero'
$ws.Cells.Item(105, 9).Value = 'Primera'
$ws.Cells.Item(105, 10).Value = 33000
$ws.Cells.Item(105, 11).Value = 300
$ws.Cells.Item(105, 12).Value = 350
$ws.Cells.Item(105, 13).Value = 327
$ws.Cells.Item(105, 14).Value = '$/unidad'
$ws.Cells.Item(105, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(105, 16).Value = 327
$ws.Cells.Item(105, 17).Value = 1
$ws.Cells.Item(105, 18).Value = 'Hortaliza'

# Row 106
$ws.Cells.Item(106, 1).Value = 7
$ws.Cells.Item(106, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(106, 3).Value = 'Ñuble'
$ws.Cells.Item(106, 4).Value = 44203
$ws.Cells.Item(106, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(106, 5).Value = 16
$ws.Cells.Item(106, 6).Value = 100112024
$ws.Cells.Item(106, 7).Value = 'Choclo'
$ws.Cells.Item(106, 8).Value = 'Choclero'
$ws.Cells.Item(106, 9).Value = 'Segunda'
$ws.Cells.Item(106, 10).Value = 12000
$ws.Cells.Item(106, 11).Value = 250
$ws.Cells.Item(106, 12).Value = 250
$ws.Cells.Item(106, 13).Value = 250
$ws.Cells.Item(106, 14).Value = '$/unidad'
$ws.Cells.Item(106, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(106, 16).Value = 250
$ws.Cells.Item(106, 17).Value = 1
$ws.Cells.Item(106, 18).Value = 'Hortaliza'

# Row 107
$ws.Cells.Item(107, 1).Value = 7
$ws.Cells.Item(107, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(107, 3).Value = 'Ñuble'
$ws.Cells.Item(107, 4).Value = 44519
$ws.Cells.Item(107, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 5).Value = 16
$ws.Cells.Item(107, 6).Value = 100112024
$ws.Cells.Item(107, 7).Value = 'Choclo'
$ws.Cells.Item(107, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(107, 9).Value = 'Primera'
$ws.Cells.Item(107, 10).Value = 100
$ws.Cells.Item(107, 11).Value = 16000
$ws.Cells.Item(107, 12).Value = 17000
$ws.Cells.Item(107, 13).Value = 16500
$ws.Cells.Item(107, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(107, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(107, 16).Value = 275
$ws.Cells.Item(107, 17).Value = 60
$ws.Cells.Item(107, 18).Value = 'Hortaliza'

# Row 108
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(108, 3).Value = 'Ñuble'
$ws.Cells.Item(108, 4).Value = 44545
$ws.Cells.Item(108, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112024
$ws.Cells.Item(108, 7).Value = 'Choclo'
$ws.Cells.Item(108, 8).Value = 'Choclero'
$ws.Cells.Item(108, 9).Value = 'Primera'
$ws.Cells.Item(108, 10).Value = 10000
$ws.Cells.Item(108, 11).Value = 350
$ws.Cells.Item(108, 12).Value = 400
$ws.Cells.Item(108, 13).Value = 375
$ws.Cells.Item(108, 14).Value = '$/unidad'
$ws.Cells.Item(108, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(108, 16).Value = 375
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = 'Hortaliza'

# Row 109
$ws.Cells.Item(109, 1).Value = 7
$ws.Cells.Item(109, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(109, 3).Value = 'Ñuble'
$ws.Cells.Item(109, 4).Value = 44253
$ws.Cells.Item(109, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(109, 5).Value = 16
$ws.Cells.Item(109, 6).Value = 100112024
$ws.Cells.Item(109, 7).Value = 'Choclo'
$ws.Cells.Item(109, 8).Value = 'Choclero'
$ws.Cells.Item(109, 9).Value = 'Primera'
$ws.Cells.Item(109, 10).Value = 47000
$ws.Cells.Item(109, 11).Value = 250
$ws.Cells.Item(109, 12).Value = 280
$ws.Cells.Item(109, 13).Value = 266
$ws.Cells.Item(109, 14).Value = '$/unidad'
$ws.Cells.Item(109, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(109, 16).Value = 266
$ws.Cells.Item(109, 17).Value = 1
$ws.Cells.Item(109, 18).Value = 'Hortaliza'

# Row 110
$ws.Cells.Item(110, 1).Value = 7
$ws.Cells.Item(110, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(110, 3).Value = 'Ñuble'
$ws.Cells.Item(110, 4).Value = 44581
$ws.Cells.Item(110, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(110, 5).Value = 16
$ws.Cells.Item(110, 6).Value = 100112024
$ws.Cells.Item(110, 7).Value = 'Choclo'
$ws.Cells.Item(110, 8).Value = 'Choclero'
$ws.Cells.Item(110, 9).Value = 'Primera'
$ws.Cells.Item(110, 10).Value = 10000
$ws.Cells.Item(110, 11).Value = 250
$ws.Cells.Item(110, 12).Value = 300
$ws.Cells.Item(110, 13).Value = 275
$ws.Cells.Item(110, 14).Value = '$/unidad'
$ws.Cells.Item(110, 15).Value = 'Región del Maule'
$ws.Cells.Item(110, 16).Value = 275
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = 'Hortaliza'

# Row 111
$ws.Cells.Item(111, 1).Value = 7
$ws.Cells.Item(111, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(111, 3).Value = 'Ñuble'
$ws.Cells.Item(111, 4).Value = 44581
$ws.Cells.Item(111, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = 100112024
$ws.Cells.Item(111, 7).Value = 'Choclo'
$ws.Cells.Item(111, 8).Value = 'Choclero'
$ws.Cells.Item(111, 9).Value = 'Segunda'
$ws.Cells.Item(111, 10).Value = 12000
$ws.Cells.Item(111, 11).Value = 150
$ws.Cells.Item(111, 12).Value = 200
$ws.Cells.Item(111, 13).Value = 175
$ws.Cells.Item(111, 14).Value = '$/unidad'
$ws.Cells.Item(111, 15).Value = 'Región del Maule'
$ws.Cells.Item(111, 16).Value = 175
$ws.Cells.Item(111, 17).Value = 1
$ws.Cells.Item(111, 18).Value = 'Hortaliza'

# Row 112
$ws.Cells.Item(112, 1).Value = 7
$ws.Cells.Item(112, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(112, 3).Value = 'Ñuble'
$ws.Cells.Item(112, 4).Value = 44567
$ws.Cells.Item(112, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 6).Value = 100112024
$ws.Cells.Item(112, 7).Value = 'Choclo'
$ws.Cells.Item(112, 8).Value = 'Choclero'
$ws.Cells.Item(112, 9).Value = 'Primera'
$ws.Cells.Item(112, 10).Value = 10000
$ws.Cells.Item(112, 11).Value = 300
$ws.Cells.Item(112, 12).Value = 350
$ws.Cells.Item(112, 13).Value = 325
$ws.Cells.Item(112, 14).Value = '$/unidad'
$ws.Cells.Item(112, 15).Value = 'Región del Maule'
$ws.Cells.Item(112, 16).Value = 325
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = 'Hortaliza'

# Row 113
$ws.Cells.Item(113, 1).Value = 7
$ws.Cells.Item(113, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(113, 3).Value = 'Ñuble'
$ws.Cells.Item(113, 4).Value = 44266
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = 100112024
$ws.Cells.Item(113, 7).Value = 'Choclo'
$ws.Cells.Item(113, 8).Value = 'Choclero'
$ws.Cells.Item(113, 9).Value = 'Primera'
$ws.Cells.Item(113, 10).Value = 48000
$ws.Cells.Item(113, 11).Value = 250
$ws.Cells.Item(113, 12).Value = 280
$ws.Cells.Item(113, 13).Value = 264
$ws.Cells.Item(113, 14).Value = '$/unidad'
$ws.Cells.Item(113, 15).Value = 'Región del Maule'
$ws.Cells.Item(113, 16).Value = 264
$ws.Cells.Item(113, 17).Value = 1
$ws.Cells.Item(113, 18).Value = 'Hortaliza'

# Row 114
$ws.Cells.Item(114, 1).Value = 7
$ws.Cells.Item(114, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(114, 3).Value = 'Ñuble'
$ws.Cells.Item(114, 4).Value = 44202
$ws.Cells.Item(114, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(114, 5).Value = 16
$ws.Cells.Item(114, 6).Value = 100112024
$ws.Cells.Item(114, 7).Value = 'Choclo'
$ws.Cells.Item(114, 8).Value = 'Choclero'
$ws.Cells.Item(114, 9).Value = 'Primera'
$ws.Cells.Item(114, 10).Value = 53000
$ws.Cells.Item(114, 11).Value = 300
$ws.Cells.Item(114, 12).Value = 350
$ws.Cells.Item(114, 13).Value = 317
$ws.Cells.Item(114, 14).Value = '$/unidad'
$ws.Cells.Item(114, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(114, 16).Value = 317
$ws.Cells.Item(114, 17).Value = 1
$ws.Cells.Item(114, 18).Value = 'Hortaliza'

# Row 115
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(115, 3).Value = 'Ñuble'
$ws.Cells.Item(115, 4).Value = 44202
$ws.Cells.Item(115, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = 100112024
$ws.Cells.Item(115, 7).Value = 'Choclo'
$ws.Cells.Item(115, 8).Value = 'Choclero'
$ws.Cells.Item(115, 9).Value = 'Segunda'
$ws.Cells.Item(115, 10).Value = 15000
$ws.Cells.Item(115, 11).Value = 250
$ws.Cells.Item(115, 12).Value = 250
$ws.Cells.Item(115, 13).Value = 250
$ws.Cells.Item(115, 14).Value = '$/unidad'
$ws.Cells.Item(115, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(115, 16).Value = 250
$ws.Cells.Item(115, 17).Value = 1
$ws.Cells.Item(115, 18).Value = 'Hortaliza'

# Row 116
$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(116, 3).Value = 'Ñuble'
$ws.Cells.Item(116, 4).Value = 44202
$ws.Cells.Item(116, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100112024
$ws.Cells.Item(116, 7).Value = 'Choclo'
$ws.Cells.Item(116, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(116, 9).Value = 'Primera'
$ws.Cells.Item(116, 10).Value = 35000
$ws.Cells.Item(116, 11).Value = 270
$ws.Cells.Item(116, 12).Value = 300
$ws.Cells.Item(116, 13).Value = 283
$ws.Cells.Item(116, 14).Value = '$/unidad'
$ws.Cells.Item(116, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(116, 16).Value = 283
$ws.Cells.Item(116, 17).Value = 1
$ws.Cells.Item(116, 18).Value = 'Hortaliza'

# Row 117
$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(117, 3).Value = 'Ñuble'
$ws.Cells.Item(117, 4).Value = 44249
$ws.Cells.Item(117, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100112024
$ws.Cells.Item(117, 7).Value = 'Choclo'
$ws.Cells.Item(117, 8).Value = 'Choclero'
$ws.Cells.Item(117, 9).Value = 'Primera'
$ws.Cells.Item(117, 10).Value = 12000
$ws.Cells.Item(117, 11).Value = 340
$ws.Cells.Item(117, 12).Value = 350
$ws.Cells.Item(117, 13).Value = 345
$ws.Cells.Item(117, 14).Value = '$/unidad'
$ws.Cells.Item(117, 15).Value = 'Región del Maule'
$ws.Cells.Item(117, 16).Value = 345
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = 'Hortaliza'

# Row 118
$ws.Cells.Item(118, 1).Value = 7
$ws.Cells.Item(118, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(118, 3).Value = 'Ñuble'
$ws.Cells.Item(118, 4).Value = 44225
$ws.Cells.Item(118, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(118, 5).Value = 16
$ws.Cells.Item(118, 6).Value = 100112024
$ws.Cells.Item(118, 7).Value = 'Choclo'
$ws.Cells.Item(118, 8).Value = 'Choclero'
$ws.Cells.Item(118, 9).Value = 'Primera'
$ws.Cells.Item(118, 10).Value = 53000
$ws.Cells.Item(118, 11).Value = 240
$ws.Cells.Item(118, 12).Value = 280
$ws.Cells.Item(118, 13).Value = 259
$ws.Cells.Item(118, 14).Value = '$/unidad'
$ws.Cells.Item(118, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(118, 16).Value = 259
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = 'Hortaliza'

# Row 119
$ws.Cells.Item(119, 1).Value = 7
$ws.Cells.Item(119, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(119, 3).Value = 'Ñuble'
$ws.Cells.Item(119, 4).Value = 44225
$ws.Cells.Item(119, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 6).Value = 100112024
$ws.Cells.Item(119, 7).Value = 'Choclo'
$ws.Cells.Item(119, 8).Value = 'Choclero'
$ws.Cells.Item(119, 9).Value = 'Segunda'
$ws.Cells.Item(119, 10).Value = 18000
$ws.Cells.Item(119, 11).Value = 200
$ws.Cells.Item(119, 12).Value = 200
$ws.Cells.Item(119, 13).Value = 200
$ws.Cells.Item(119, 14).Value = '$/unidad'
$ws.Cells.Item(119, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(119, 16).Value = 200
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = 'Hortaliza'

# Row 120
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(120, 3).Value = 'Ñuble'
$ws.Cells.Item(120, 4).Value = 44259
$ws.Cells.Item(120, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112024
$ws.Cells.Item(120, 7).Value = 'Choclo'
$ws.Cells.Item(120, 8).Value = 'Choclero'
$ws.Cells.Item(120, 9).Value = 'Primera'
$ws.Cells.Item(120, 10).Value = 54000
$ws.Cells.Item(120, 11).Value = 250
$ws.Cells.Item(120, 12).Value = 280
$ws.Cells.Item(120, 13).Value = 260
$ws.Cells.Item(120, 14).Value = '$/unidad'
$ws.Cells.Item(120, 15).Value = 'Región del Maule'
$ws.Cells.Item(120, 16).Value = 260
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = 'Hortaliza'

# Row 121
$ws.Cells.Item(121, 1).Value = 7
$ws.Cells.Item(121, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(121, 3).Value = 'Ñuble'
$ws.Cells.Item(121, 4).Value = 44328
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 5).Value = 16
$ws.Cells.Item(121, 6).Value = 100112024
$ws.Cells.Item(121, 7).Value = 'Choclo'
$ws.Cells.Item(121, 8).Value = 'Choclero'
$ws.Cells.Item(121, 9).Value = 'Primera'
$ws.Cells.Item(121, 10).Value = 3000
$ws.Cells.Item(121, 11).Value = 180
$ws.Cells.Item(121, 12).Value = 200
$ws.Cells.Item(121, 13).Value = 190
$ws.Cells.Item(121, 14).Value = '$/unidad'
$ws.Cells.Item(121, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(121, 16).Value = 190
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = 'Hortaliza'

# Row 122
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(122, 3).Value = 'Ñuble'
$ws.Cells.Item(122, 4).Value = 44561
$ws.Cells.Item(122, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = 100112024
$ws.Cells.Item(122, 7).Value = 'Choclo'
$ws.Cells.Item(122, 8).Value = 'Choclero'
$ws.Cells.Item(122, 9).Value = 'Segunda'
$ws.Cells.Item(122, 10).Value = 12000
$ws.Cells.Item(122, 11).Value = 200
$ws.Cells.Item(122, 12).Value = 250
$ws.Cells.Item(122, 13).Value = 225
$ws.Cells.Item(122, 14).Value = '$/unidad'
$ws.Cells.Item(122, 15).Value = 'Región del Maule'
$ws.Cells.Item(122, 16).Value = 225
$ws.Cells.Item(122, 17).Value = 1
$ws.Cells.Item(122, 18).Value = 'Hortaliza'

# Row 123
$ws.Cells.Item(123, 1).Value = 7
$ws.Cells.Item(123, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(123, 3).Value = 'Ñuble'
$ws.Cells.Item(123, 4).Value = 44526
$ws.Cells.Item(123, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(123, 5).Value = 16
$ws.Cells.Item(123, 6).Value = 100112024
$ws.Cells.Item(123, 7).Value = 'Choclo'
$ws.Cells.Item(123, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(123, 9).Value = 'Primera'
$ws.Cells.Item(123, 10).Value = 60
$ws.Cells.Item(123, 11).Value = 16000
$ws.Cells.Item(123, 12).Value = 17000
$ws.Cells.Item(123, 13).Value = 16500
$ws.Cells.Item(123, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(123, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(123, 16).Value = 275
$ws.Cells.Item(123, 17).Value = 60
$ws.Cells.Item(123, 18).Value = 'Hortaliza'

# Row 124
$ws.Cells.Item(124, 1).Value = 7
$ws.Cells.Item(124, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(124, 3).Value = 'Ñuble'
$ws.Cells.Item(124, 4).Value = 44250
$ws.Cells.Item(124, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 5).Value = 16
$ws.Cells.Item(124, 6).Value = 100112024
$ws.Cells.Item(124, 7).Value = 'Choclo'
$ws.Cells.Item(124, 8).Value = 'Choclero'
$ws.Cells.Item(124, 9).Value = 'Primera'
$ws.Cells.Item(124, 10).Value = 16000
$ws.Cells.Item(124, 11).Value = 350
$ws.Cells.Item(124, 12).Value = 360
$ws.Cells.Item(124, 13).Value = 355
$ws.Cells.Item(124, 14).Value = '$/unidad'
$ws.Cells.Item(124, 15).Value = 'Región del Maule'
$ws.Cells.Item(124, 16).Value = 355
$ws.Cells.Item(124, 17).Value = 1
$ws.Cells.Item(124, 18).Value = 'Hortaliza'

# Row 125
$ws.Cells.Item(125, 1).Value = 7
$ws.Cells.Item(125, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(125, 3).Value = 'Ñuble'
$ws.Cells.Item(125, 4).Value = 44285
$ws.Cells.Item(125, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 5).Value = 16
$ws.Cells.Item(125, 6).Value = 100112024
$ws.Cells.Item(125, 7).Value = 'Choclo'
$ws.Cells.Item(125, 8).Value = 'Choclero'
$ws.Cells.Item(125, 9).Value = 'Primera'
$ws.Cells.Item(125, 10).Value = 12000
$ws.Cells.Item(125, 11).Value = 180
$ws.Cells.Item(125, 12).Value = 200
$ws.Cells.Item(125, 13).Value = 190
$ws.Cells.Item(125, 14).Value = '$/unidad'
$ws.Cells.Item(125, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(125, 16).Value = 190
$ws.Cells.Item(125, 17).Value = 1
$ws.Cells.Item(125, 18).Value = 'Hortaliza'

# Row 126
$ws.Cells.Item(126, 1).Value = 7
$ws.Cells.Item(126, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(126, 3).Value = 'Ñuble'
$ws.Cells.Item(126, 4).Value = 44264
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 5).Value = 16
$ws.Cells.Item(126, 6).Value = 100112024
$ws.Cells.Item(126, 7).Value = 'Choclo'
$ws.Cells.Item(126, 8).Value = 'Choclero'
$ws.Cells.Item(126, 9).Value = 'Primera'
$ws.Cells.Item(126, 10).Value = 50000
$ws.Cells.Item(126, 11).Value = 270
$ws.Cells.Item(126, 12).Value = 300
$ws.Cells.Item(126, 13).Value = 291
$ws.Cells.Item(126, 14).Value = '$/unidad'
$ws.Cells.Item(126, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(126, 16).Value = 291
$ws.Cells.Item(126, 17).Value = 1
$ws.Cells.Item(126, 18).Value = 'Hortaliza'

# Row 127
$ws.Cells.Item(127, 1).Value = 7
$ws.Cells.Item(127, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(127, 3).Value = 'Ñuble'
$ws.Cells.Item(127, 4).Value = 44264
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 5).Value = 16
$ws.Cells.Item(127, 6).Value = 100112024
$ws.Cells.Item(127, 7).Value = 'Choclo'
$ws.Cells.Item(127, 8).Value = 'Choclero'
$ws.Cells.Item(127, 9).Value = 'Segunda'
$ws.Cells.Item(127, 10).Value = 33000
$ws.Cells.Item(127, 11).Value = 220
$ws.Cells.Item(127, 12).Value = 250
$ws.Cells.Item(127, 13).Value = 236
$ws.Cells.Item(127, 14).Value = '$/unidad'
$ws.Cells.Item(127, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(127, 16).Value = 236
$ws.Cells.Item(127, 17).Value = 1
$ws.Cells.Item(127, 18).Value = 'Hortaliza'

# Row 128
$ws.Cells.Item(128, 1).Value = 7
$ws.Cells.Item(128, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(128, 3).Value = 'Ñuble'
$ws.Cells.Item(128, 4).Value = 44533
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(128, 5).Value = 16
$ws.Cells.Item(128, 6).Value = 100112024
$ws.Cells.Item(128, 7).Value = 'Choclo'
$ws.Cells.Item(128, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(128, 9).Value = 'Primera'
$ws.Cells.Item(128, 10).Value = 60
$ws.Cells.Item(128, 11).Value = 16000
$ws.Cells.Item(128, 12).Value = 17000
$ws.Cells.Item(128, 13).Value = 16500
$ws.Cells.Item(128, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(128, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(128, 16).Value = 275
$ws.Cells.Item(128, 17).Value = 60
$ws.Cells.Item(128, 18).Value = 'Hortaliza'

# Row 129
$ws.Cells.Item(129, 1).Value = 7
$ws.Cells.Item(129, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(129, 3).Value = 'Ñuble'
$ws.Cells.Item(129, 4).Value = 44221
$ws.Cells.Item(129, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(129, 5).Value = 16
$ws.Cells.Item(129, 6).Value = 100112024
$ws.Cells.Item(129, 7).Value = 'Choclo'
$ws.Cells.Item(129, 8).Value = 'Choclero'
$ws.Cells.Item(129, 9).Value = 'Primera'
$ws.Cells.Item(129, 10).Value = 50000
$ws.Cells.Item(129, 11).Value = 250
$ws.Cells.Item(129, 12).Value = 280
$ws.Cells.Item(129, 13).Value = 271
$ws.Cells.Item(129, 14).Value = '$/unidad'
$ws.Cells.Item(129, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(129, 16).Value = 271
$ws.Cells.Item(129, 17).Value = 1
$ws.Cells.Item(129, 18).Value = 'Hortaliza'

# Row 130
$ws.Cells.Item(130, 1).Value = 7
$ws.Cells.Item(130, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(130, 3).Value = 'Ñuble'
$ws.Cells.Item(130, 4).Value = 44523
$ws.Cells.Item(130, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(130, 5).Value = 16
$ws.Cells.Item(130, 6).Value = 100112024
$ws.Cells.Item(130, 7).Value = 'Choclo'
$ws.Cells.Item(130, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(130, 9).Value = 'Primera'
$ws.Cells.Item(130, 10).Value = 80
$ws.Cells.Item(130, 11).Value = 16000
$ws.Cells.Item(130, 12).Value = 17000
$ws.Cells.Item(130, 13).Value = 16500
$ws.Cells.Item(130, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(130, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(130, 16).Value = 275
$ws.Cells.Item(130, 17).Value = 60
$ws.Cells.Item(130, 18).Value = 'Hortaliza'

# Row 131
$ws.Cells.Item(131, 1).Value = 7
$ws.Cells.Item(131, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(131, 3).Value = 'Ñuble'
$ws.Cells.Item(131, 4).Value = 44316
$ws.Cells.Item(131, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(131, 5).Value = 16
$ws.Cells.Item(131, 6).Value = 100112024
$ws.Cells.Item(131, 7).Value = 'Choclo'
$ws.Cells.Item(131, 8).Value = 'Choclero'
$ws.Cells.Item(131, 9).Value = 'Primera'
$ws.Cells.Item(131, 10).Value = 6000
$ws.Cells.Item(131, 11).Value = 200
$ws.Cells.Item(131, 12).Value = 220
$ws.Cells.Item(131, 13).Value = 210
$ws.Cells.Item(131, 14).Value = '$/unidad'
$ws.Cells.Item(131, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(131, 16).Value = 210
$ws.Cells.Item(131, 17).Value = 1
$ws.Cells.Item(131, 18).Value = 'Hortaliza'

# Row 132
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(132, 3).Value = 'Ñuble'
$ws.Cells.Item(132, 4).Value = 44186
$ws.Cells.Item(132, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112024
$ws.Cells.Item(132, 7).Value = 'Choclo'
$ws.Cells.Item(132, 8).Value = 'Choclero'
$ws.Cells.Item(132, 9).Value = 'Primera'
$ws.Cells.Item(132, 10).Value = 16000
$ws.Cells.Item(132, 11).Value = 280
$ws.Cells.Item(132, 12).Value = 300
$ws.Cells.Item(132, 13).Value = 290
$ws.Cells.Item(132, 14).Value = '$/unidad'
$ws.Cells.Item(132, 15).Value = 'Región del Maule'
$ws.Cells.Item(132, 16).Value = 290
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = 'Hortaliza'

# Row 133
$ws.Cells.Item(133, 1).Value = 7
$ws.Cells.Item(133, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(133, 3).Value = 'Ñuble'
$ws.Cells.Item(133, 4).Value = 44278
$ws.Cells.Item(133, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(133, 5).Value = 16
$ws.Cells.Item(133, 6).Value = 100112024
$ws.Cells.Item(133, 7).Value = 'Choclo'
$ws.Cells.Item(133, 8).Value = 'Choclero'
$ws.Cells.Item(133, 9).Value = 'Primera'
$ws.Cells.Item(133, 10).Value = 3000
$ws.Cells.Item(133, 11).Value = 250
$ws.Cells.Item(133, 12).Value = 260
$ws.Cells.Item(133, 13).Value = 255
$ws.Cells.Item(133, 14).Value = '$/unidad'
$ws.Cells.Item(133, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(133, 16).Value = 255
$ws.Cells.Item(133, 17).Value = 1
$ws.Cells.Item(133, 18).Value = 'Hortaliza'

# Row 134
$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(134, 3).Value = 'Ñuble'
$ws.Cells.Item(134, 4).Value = 44312
$ws.Cells.Item(134, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = 100112024
$ws.Cells.Item(134, 7).Value = 'Choclo'
$ws.Cells.Item(134, 8).Value = 'Choclero'
$ws.Cells.Item(134, 9).Value = 'Primera'
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 230
$ws.Cells.Item(134, 12).Value = 250
$ws.Cells.Item(134, 13).Value = 240
$ws.Cells.Item(134, 14).Value = '$/unidad'
$ws.Cells.Item(134, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(134, 16).Value = 240
$ws.Cells.Item(134, 17).Value = 1
$ws.Cells.Item(134, 18).Value = 'Hortaliza'

# Row 135
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(135, 3).Value = 'Ñuble'
$ws.Cells.Item(135, 4).Value = 44300
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112024
$ws.Cells.Item(135, 7).Value = 'Choclo'
$ws.Cells.Item(135, 8).Value = 'Choclero'
$ws.Cells.Item(135, 9).Value = 'Primera'
$ws.Cells.Item(135, 10).Value = 8000
$ws.Cells.Item(135, 11).Value = 180
$ws.Cells.Item(135, 12).Value = 200
$ws.Cells.Item(135, 13).Value = 190
$ws.Cells.Item(135, 14).Value = '$/unidad'
$ws.Cells.Item(135, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(135, 16).Value = 190
$ws.Cells.Item(135, 17).Value = 1
$ws.Cells.Item(135, 18).Value = 'Hortaliza'

# Row 136
$ws.Cells.Item(136, 1).Value = 7
$ws.Cells.Item(136, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(136, 3).Value = 'Ñuble'
$ws.Cells.Item(136, 4).Value = 44314
$ws.Cells.Item(136, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(136, 5).Value = 16
$ws.Cells.Item(136, 6).Value = 100112024
$ws.Cells.Item(136, 7).Value = 'Choclo'
$ws.Cells.Item(136, 8).Value = 'Choclero'
$ws.Cells.Item(136, 9).Value = 'Primera'
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 240
$ws.Cells.Item(136, 12).Value = 250
$ws.Cells.Item(136, 13).Value = 245
$ws.Cells.Item(136, 14).Value = '$/unidad'
$ws.Cells.Item(136, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(136, 16).Value = 245
$ws.Cells.Item(136, 17).Value = 1
$ws.Cells.Item(136, 18).Value = 'Hortaliza'

# Row 137
$ws.Cells.Item(137, 1).Value = 7
$ws.Cells.Item(137, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(137, 3).Value = 'Ñuble'
$ws.Cells.Item(137, 4).Value = 44260
$ws.Cells.Item(137, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(137, 5).Value = 16
$ws.Cells.Item(137, 6).Value = 100112024
$ws.Cells.Item(137, 7).Value = 'Choclo'
$ws.Cells.Item(137, 8).Value = 'Choclero'
$ws.Cells.Item(137, 9).Value = 'Primera'
$ws.Cells.Item(137, 10).Value = 90000
$ws.Cells.Item(137, 11).Value = 270
$ws.Cells.Item(137, 12).Value = 300
$ws.Cells.Item(137, 13).Value = 282
$ws.Cells.Item(137, 14).Value = '$/unidad'
$ws.Cells.Item(137, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(137, 16).Value = 282
$ws.Cells.Item(137, 17).Value = 1
$ws.Cells.Item(137, 18).Value = 'Hortaliza'

# Row 138
$ws.Cells.Item(138, 1).Value = 7
$ws.Cells.Item(138, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(138, 3).Value = 'Ñuble'
$ws.Cells.Item(138, 4).Value = 44260
$ws.Cells.Item(138, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(138, 5).Value = 16
$ws.Cells.Item(138, 6).Value = 100112024
$ws.Cells.Item(138, 7).Value = 'Choclo'
$ws.Cells.Item(138, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(138, 9).Value = 'Primera'
$ws.Cells.Item(138, 10).Value = 33000
$ws.Cells.Item(138, 11).Value = 150
$ws.Cells.Item(138, 12).Value = 180
$ws.Cells.Item(138, 13).Value = 164
$ws.Cells.Item(138, 14).Value = '$/unidad'
$ws.Cells.Item(138, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(138, 16).Value = 164
$ws.Cells.Item(138, 17).Value = 1
$ws.Cells.Item(138, 18).Value = 'Hortaliza'

# Row 139
$ws.Cells.Item(139, 1).Value = 7
$ws.Cells.Item(139, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(139, 3).Value = 'Ñuble'
$ws.Cells.Item(139, 4).Value = 44585
$ws.Cells.Item(139, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(139, 5).Value = 16
$ws.Cells.Item(139, 6).Value = 100112024
$ws.Cells.Item(139, 7).Value = 'Choclo'
$ws.Cells.Item(139, 8).Value = 'Choclero'
$ws.Cells.Item(139, 9).Value = 'Primera'
$ws.Cells.Item(139, 10).Value = 20000
$ws.Cells.Item(139, 11).Value = 250
$ws.Cells.Item(139, 12).Value = 300
$ws.Cells.Item(139, 13).Value = 275
$ws.Cells.Item(139, 14).Value = '$/unidad'
$ws.Cells.Item(139, 15).Value = 'Región del Maule'
$ws.Cells.Item(139, 16).Value = 275
$ws.Cells.Item(139, 17).Value = 1
$ws.Cells.Item(139, 18).Value = 'Hortaliza'

# Row 140
$ws.Cells.Item(140, 1).Value = 7
$ws.Cells.Item(140, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(140, 3).Value = 'Ñuble'
$ws.Cells.Item(140, 4).Value = 44585
$ws.Cells.Item(140, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(140, 5).Value = 16
$ws.Cells.Item(140, 6).Value = 100112024
$ws.Cells.Item(140, 7).Value = 'Choclo'
$ws.Cells.Item(140, 8).Value = 'Choclero'
$ws.Cells.Item(140, 9).Value = 'Segunda'
$ws.Cells.Item(140, 10).Value = 20000
$ws.Cells.Item(140, 11).Value = 150
$ws.Cells.Item(140, 12).Value = 200
$ws.Cells.Item(140, 13).Value = 175
$ws.Cells.Item(140, 14).Value = '$/unidad'
$ws.Cells.Item(140, 15).Value = 'Región del Maule'
$ws.Cells.Item(140, 16).Value = 175
$ws.Cells.Item(140, 17).Value = 1
$ws.Cells.Item(140, 18).Value = 'Hortaliza'

# Row 141
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(141, 3).Value = 'Ñuble'
$ws.Cells.Item(141, 4).Value = 44560
$ws.Cells.Item(141, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112024
$ws.Cells.Item(141, 7).Value = 'Choclo'
$ws.Cells.Item(141, 8).Value = 'Choclero'
$ws.Cells.Item(141, 9).Value = 'Primera'
$ws.Cells.Item(141, 10).Value = 14000
$ws.Cells.Item(141, 11).Value = 300
$ws.Cells.Item(141, 12).Value = 350
$ws.Cells.Item(141, 13).Value = 325
$ws.Cells.Item(141, 14).Value = '$/unidad'
$ws.Cells.Item(141, 15).Value = 'Región del Maule'
$ws.Cells.Item(141, 16).Value = 325
$ws.Cells.Item(141, 17).Value = 1
$ws.Cells.Item(141, 18).Value = 'Hortaliza'

# Row 142
$ws.Cells.Item(142, 1).Value = 7
$ws.Cells.Item(142, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(142, 3).Value = 'Ñuble'
$ws.Cells.Item(142, 4).Value = 44560
$ws.Cells.Item(142, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(142, 5).Value = 16
$ws.Cells.Item(142, 6).Value = 100112024
$ws.Cells.Item(142, 7).Value = 'Choclo'
$ws.Cells.Item(142, 8).Value = 'Choclero'
$ws.Cells.Item(142, 9).Value = 'Segunda'
$ws.Cells.Item(142, 10).Value = 10000
$ws.Cells.Item(142, 11).Value = 200
$ws.Cells.Item(142, 12).Value = 250
$ws.Cells.Item(142, 13).Value = 225
$ws.Cells.Item(142, 14).Value = '$/unidad'
$ws.Cells.Item(142, 15).Value = 'Región del Maule'
$ws.Cells.Item(142, 16).Value = 225
$ws.Cells.Item(142, 17).Value = 1
$ws.Cells.Item(142, 18).Value = 'Hortaliza'

# Row 143
$ws.Cells.Item(143, 1).Value = 7
$ws.Cells.Item(143, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(143, 3).Value = 'Ñuble'
$ws.Cells.Item(143, 4).Value = 44272
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(143, 5).Value = 16
$ws.Cells.Item(143, 6).Value = 100112024
$ws.Cells.Item(143, 7).Value = 'Choclo'
$ws.Cells.Item(143, 8).Value = 'Choclero'
$ws.Cells.Item(143, 9).Value = 'Primera'
$ws.Cells.Item(143, 10).Value = 65000
$ws.Cells.Item(143, 11).Value = 240
$ws.Cells.Item(143, 12).Value = 260
$ws.Cells.Item(143, 13).Value = 249
$ws.Cells.Item(143, 14).Value = '$/unidad'
$ws.Cells.Item(143, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(143, 16).Value = 249
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = 'Hortaliza'

# Row 144
$ws.Cells.Item(144, 1).Value = 7
$ws.Cells.Item(144, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(144, 3).Value = 'Ñuble'
$ws.Cells.Item(144, 4).Value = 44272
$ws.Cells.Item(144, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(144, 5).Value = 16
$ws.Cells.Item(144, 6).Value = 100112024
$ws.Cells.Item(144, 7).Value = 'Choclo'
$ws.Cells.Item(144, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(144, 9).Value = 'Primera'
$ws.Cells.Item(144, 10).Value = 70000
$ws.Cells.Item(144, 11).Value = 130
$ws.Cells.Item(144, 12).Value = 150
$ws.Cells.Item(144, 13).Value = 143
$ws.Cells.Item(144, 14).Value = '$/unidad'
$ws.Cells.Item(144, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(144, 16).Value = 143
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = 'Hortaliza'

# Row 145
$ws.Cells.Item(145, 1).Value = 7
$ws.Cells.Item(145, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(145, 3).Value = 'Ñuble'
$ws.Cells.Item(145, 4).Value = 44162
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(145, 5).Value = 16
$ws.Cells.Item(145, 6).Value = 100112024
$ws.Cells.Item(145, 7).Value = 'Choclo'
$ws.Cells.Item(145, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(145, 9).Value = 'Primera'
$ws.Cells.Item(145, 10).Value = 40
$ws.Cells.Item(145, 11).Value = 23000
$ws.Cells.Item(145, 12).Value = 24000
$ws.Cells.Item(145, 13).Value = 23450
$ws.Cells.Item(145, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(145, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(145, 16).Value = 335
$ws.Cells.Item(145, 17).Value = 70
$ws.Cells.Item(145, 18).Value = 'Hortaliza'

# Row 146
$ws.Cells.Item(146, 1).Value = 7
$ws.Cells.Item(146, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(146, 3).Value = 'Ñuble'
$ws.Cells.Item(146, 4).Value = 44529
$ws.Cells.Item(146, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(146, 5).Value = 16
$ws.Cells.Item(146, 6).Value = 100112024
$ws.Cells.Item(146, 7).Value = 'Choclo'
$ws.Cells.Item(146, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(146, 9).Value = 'Primera'
$ws.Cells.Item(146, 10).Value = 80
$ws.Cells.Item(146, 11).Value = 16000
$ws.Cells.Item(146, 12).Value = 17000
$ws.Cells.Item(146, 13).Value = 16500
$ws.Cells.Item(146, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(146, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(146, 16).Value = 275
$ws.Cells.Item(146, 17).Value = 60
$ws.Cells.Item(146, 18).Value = 'Hortaliza'

# Row 147
$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(147, 3).Value = 'Ñuble'
$ws.Cells.Item(147, 4).Value = 44189
$ws.Cells.Item(147, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = 100112024
$ws.Cells.Item(147, 7).Value = 'Choclo'
$ws.Cells.Item(147, 8).Value = 'Choclero'
$ws.Cells.Item(147, 9).Value = 'Primera'
$ws.Cells.Item(147, 10).Value = 12000
$ws.Cells.Item(147, 11).Value = 380
$ws.Cells.Item(147, 12).Value = 400
$ws.Cells.Item(147, 13).Value = 390
$ws.Cells.Item(147, 14).Value = '$/unidad'
$ws.Cells.Item(147, 15).Value = 'Región del Maule'
$ws.Cells.Item(147, 16).Value = 390
$ws.Cells.Item(147, 17).Value = 1
$ws.Cells.Item(147, 18).Value = 'Hortaliza'

# Row 148
$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(148, 3).Value = 'Ñuble'
$ws.Cells.Item(148, 4).Value = 44274
$ws.Cells.Item(148, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112024
$ws.Cells.Item(148, 7).Value = 'Choclo'
$ws.Cells.Item(148, 8).Value = 'Choclero'
$ws.Cells.Item(148, 9).Value = 'Primera'
$ws.Cells.Item(148, 10).Value = 47000
$ws.Cells.Item(148, 11).Value = 250
$ws.Cells.Item(148, 12).Value = 270
$ws.Cells.Item(148, 13).Value = 259
$ws.Cells.Item(148, 14).Value = '$/unidad'
$ws.Cells.Item(148, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(148, 16).Value = 259
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = 'Hortaliza'

# Row 149
$ws.Cells.Item(149, 1).Value = 7
$ws.Cells.Item(149, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(149, 3).Value = 'Ñuble'
$ws.Cells.Item(149, 4).Value = 44274
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 16
$ws.Cells.Item(149, 6).Value = 100112024
$ws.Cells.Item(149, 7).Value = 'Choclo'
$ws.Cells.Item(149, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(149, 9).Value = 'Primera'
$ws.Cells.Item(149, 10).Value = 33000
$ws.Cells.Item(149, 11).Value = 130
$ws.Cells.Item(149, 12).Value = 150
$ws.Cells.Item(149, 13).Value = 141
$ws.Cells.Item(149, 14).Value = '$/unidad'
$ws.Cells.Item(149, 15).Value = 'Región del Maule'
$ws.Cells.Item(149, 16).Value = 141
$ws.Cells.Item(149, 17).Value = 1
$ws.Cells.Item(149, 18).Value = 'Hortaliza'

# Row 150
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(150, 3).Value = 'Ñuble'
$ws.Cells.Item(150, 4).Value = 44554
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112024
$ws.Cells.Item(150, 7).Value = 'Choclo'
$ws.Cells.Item(150, 8).Value = 'Choclero'
$ws.Cells.Item(150, 9).Value = 'Primera'
$ws.Cells.Item(150, 10).Value = 14000
$ws.Cells.Item(150, 11).Value = 300
$ws.Cells.Item(150, 12).Value = 350
$ws.Cells.Item(150, 13).Value = 325
$ws.Cells.Item(150, 14).Value = '$/unidad'
$ws.Cells.Item(150, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(150, 16).Value = 325
$ws.Cells.Item(150, 17).Value = 1
$ws.Cells.Item(150, 18).Value = 'Hortaliza'

# Row 151
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(151, 3).Value = 'Ñuble'
$ws.Cells.Item(151, 4).Value = 44554
$ws.Cells.Item(151, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = 100112024
$ws.Cells.Item(151, 7).Value = 'Choclo'
$ws.Cells.Item(151, 8).Value = 'Choclero'
$ws.Cells.Item(151, 9).Value = 'Segunda'
$ws.Cells.Item(151, 10).Value = 5000
$ws.Cells.Item(151, 11).Value = 250
$ws.Cells.Item(151, 12).Value = 250
$ws.Cells.Item(151, 13).Value = 250
$ws.Cells.Item(151, 14).Value = '$/unidad'
$ws.Cells.Item(151, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(151, 16).Value = 250
$ws.Cells.Item(151, 17).Value = 1
$ws.Cells.Item(151, 18).Value = 'Hortaliza'
